$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") with the latest scraped values, and swap the Monero/Stacks
# rows (41/42) to reflect the new ranking order, per the GitHub Actions
# "Updated cryptos list" job.
#
# Column D holds text-formatted numbers (e.g. "67.978.26", "0.999") rather
# than real numbers, so for plain decimal-looking values we force the cell
# to Text format ("@") before assigning the string - otherwise Excel would
# auto-coerce a value like "1.00" or "614.77" into the number 1 / 614.77.

$ws.Range('D2').Value = '67.904.76'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '3.540.84'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.77'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.84'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').Value = '3.539.68'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.07'
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.427'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').Value = '4.139.26'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.10'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '3.544.99'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '67.632.25'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.38'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.73'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '447.92'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.624'
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.51'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000131'
$ws.Range('E25').Value = '  +6.41%  '
$ws.Range('D26').Value = '3.682.07'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.68'
$ws.Range('E29').Value = '  +4.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.55'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.62'
$ws.Range('E31').Value = '  -3.73%  '
$ws.Range('E32').Value = '  +7.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.96'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('D36').Value = '3.527.80'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.21'
$ws.Range('E41').Value = '  +3.44%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '175.59'
$ws.Range('E42').Value = '  -0.72%  '
$ws.Range('E43').Value = '  +2.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.43'
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.887'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.90'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.54'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.28'
$ws.Range('E49').Value = '  +4.81%  '
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.995'
$ws.Range('E51').Value = '  -3.86%  '
